$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AWAY / HOME / MATCH_RESULT rows 2-9 with the new matchups
$ws.Range("A2").Value = "Cavaliers"
$ws.Range("B2").Value = "Hornets"
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = "Pistons"
$ws.Range("B3").Value = "Wizards"
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = "Nuggets"
$ws.Range("B4").Value = "Raptors"
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = "Lakers"
$ws.Range("B5").Value = "Pelicans"
$ws.Range("C5").Value = 2

$ws.Range("A6").Value = "Nets"
$ws.Range("B6").Value = "Thunder"
$ws.Range("C6").Value = 2

$ws.Range("A7").Value = "Magic"
$ws.Range("B7").Value = "Spurs"
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = "Bucks"
$ws.Range("B8").Value = "Suns"
$ws.Range("C8").Value = 2

$ws.Range("A9").Value = "Knicks"
$ws.Range("B9").Value = "Trail Blazers"
$ws.Range("C9").Value = 2

# Remove the old row 10 (Jazz vs Kings) entirely, shrinking the used range to A1:C9
$ws.Rows(10).Delete()
